$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3953.8518
$ws.Range("I62").Value = 1466.6666
$ws.Range("J62").Value = 5943.6
$ws.Range("K62").Value = 1466.6666
$ws.Range("L62").Value = 5943.6
$ws.Range("M62").Value = -842.6666
$ws.Range("N62").Value = -7191.6

$ws.Range("H65").Value = 3953.8518
$ws.Range("I65").Value = 1466.6666
$ws.Range("J65").Value = 5943.6
$ws.Range("K65").Value = 7333.333000000001
$ws.Range("L65").Value = 29718
$ws.Range("M65").Value = -4213.333000000001
$ws.Range("N65").Value = -35958

$ws.Range("H100").Value = 3237
$ws.Range("I100").Value = 3237
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3237
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2696

$ws.Range("H107").Value = 55557040
$ws.Range("I107").Value = 58825070
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 58825070
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = -58823150
$ws.Range("N107").Value = -4340

$ws.Range("H116").Value = 2935
$ws.Range("I116").Value = 2935
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2935
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 507

$ws.Range("H121").Value = 2499.625
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2499.625
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 7498.875
$ws.Range("N121").Value = -10992.875

$ws.Range("H132").Value = 4780.2
$ws.Range("I132").Value = 3089.2222
$ws.Range("J132").Value = 19999
$ws.Range("K132").Value = 9267.6666
$ws.Range("L132").Value = 59997
$ws.Range("M132").Value = -6737.6666
$ws.Range("N132").Value = -65057

$ws.Range("H135").Value = 2580
$ws.Range("I135").Value = 1225
$ws.Range("J135").Value = 8000
$ws.Range("K135").Value = 11025
$ws.Range("L135").Value = 72000
$ws.Range("M135").Value = -8490
$ws.Range("N135").Value = -77070

$ws.Range("H137").Value = 5991.875
$ws.Range("I137").Value = 7299.8335
$ws.Range("J137").Value = 4683.9165
$ws.Range("K137").Value = 21899.5005
$ws.Range("L137").Value = 14051.7495
$ws.Range("M137").Value = -19349.5005
$ws.Range("N137").Value = -19151.7495

$ws.Range("H138").Value = 4333.3335
$ws.Range("I138").Value = 4333.3335
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 13000.0005
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -7860.000499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 22223292
$ws.Range("I2").Value = 23810564
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 23810564
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -23810451
$ws.Range("N2").Value = -1726

$ws.Range("H110").Value = 4116761.8
$ws.Range("I110").Value = 7937540.5
$ws.Range("J110").Value = 2076.7693
$ws.Range("K110").Value = 7937540.5
$ws.Range("L110").Value = 2076.7693
$ws.Range("M110").Value = -7935495.5
$ws.Range("N110").Value = -6166.7693

$ws.Range("H116").Value = 22223292
$ws.Range("I116").Value = 23810564
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 23810564
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = -23808270
$ws.Range("N116").Value = -6088

$ws.Range("H122").Value = 420010.25
$ws.Range("I122").Value = 558380.3
$ws.Range("J122").Value = 4900
$ws.Range("K122").Value = 1675140.9
$ws.Range("L122").Value = 14700
$ws.Range("M122").Value = -1672690.9
$ws.Range("N122").Value = -19600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 22223292
$ws.Range("I3").Value = 23810564
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 23810564
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -23810450
$ws.Range("N3").Value = -1728

$ws.Range("H86").Value = 1832.7778
$ws.Range("I86").Value = 1732.8
$ws.Range("J86").Value = 2332.6667
$ws.Range("K86").Value = 1732.8
$ws.Range("L86").Value = 2332.6667
$ws.Range("M86").Value = -609.8
$ws.Range("N86").Value = -4578.6667

$ws.Range("H89").Value = 1832.7778
$ws.Range("I89").Value = 1732.8
$ws.Range("J89").Value = 2332.6667
$ws.Range("K89").Value = 8664
$ws.Range("L89").Value = 11663.3335
$ws.Range("M89").Value = -3048
$ws.Range("N89").Value = -22895.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3830.6667
$ws.Range("I31").Value = 4299.9287
$ws.Range("J31").Value = 3484.8948
$ws.Range("K31").Value = 4299.9287
$ws.Range("L31").Value = 3484.8948
$ws.Range("M31").Value = -4004.9287
$ws.Range("N31").Value = -4074.8948

$ws.Range("H34").Value = 3830.6667
$ws.Range("I34").Value = 4299.9287
$ws.Range("J34").Value = 3484.8948
$ws.Range("K34").Value = 4299.9287
$ws.Range("L34").Value = 3484.8948
$ws.Range("M34").Value = -4097.9287
$ws.Range("N34").Value = -3888.8948

$ws.Range("H132").Value = 4937.5557
$ws.Range("I132").Value = 4305
$ws.Range("J132").Value = 9998
$ws.Range("K132").Value = 12915
$ws.Range("L132").Value = 29994
$ws.Range("M132").Value = -10385
$ws.Range("N132").Value = -35054

$ws.Range("H134").Value = 5749.3335
$ws.Range("I134").Value = 9999
$ws.Range("J134").Value = 3624.5
$ws.Range("K134").Value = 29997
$ws.Range("L134").Value = 10873.5
$ws.Range("M134").Value = -27462
$ws.Range("N134").Value = -15943.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 301.54544
$ws.Range("I2").Value = 11.25
$ws.Range("J2").Value = 467.42856
$ws.Range("K2").Value = 67.5
$ws.Range("L2").Value = 2804.57136
$ws.Range("M2").Value = 45.5
$ws.Range("N2").Value = -3030.57136

$ws.Range("H9").Value = 2875.25
$ws.Range("I9").Value = 1333.6666
$ws.Range("J9").Value = 7500
$ws.Range("K9").Value = 4000.9998
$ws.Range("L9").Value = 22500
$ws.Range("M9").Value = -3776.9998

$ws.Range("H14").Value = 4000
$ws.Range("I14").Value = 4000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -11827

$ws.Range("H23").Value = 202120
$ws.Range("I23").Value = 200
$ws.Range("J23").Value = 252600
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 757800
$ws.Range("M23").Value = -365
$ws.Range("N23").Value = -758270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2537.4707
$ws.Range("I102").Value = 713.8
$ws.Range("J102").Value = 5142.7144
$ws.Range("K102").Value = 713.8
$ws.Range("L102").Value = 5142.7144
$ws.Range("M102").Value = 908.2
$ws.Range("N102").Value = -8386.714400000001

$ws.Range("H107").Value = 1770.64
$ws.Range("I107").Value = 1587.7368
$ws.Range("J107").Value = 2349.8333
$ws.Range("K107").Value = 1587.7368
$ws.Range("L107").Value = 2349.8333
$ws.Range("M107").Value = 332.2632000000001
$ws.Range("N107").Value = -6189.8333

$ws.Range("H132").Value = 6007
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 6007
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 18021
$ws.Range("N132").Value = -23081
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1000.2857
$ws.Range("I113").Value = 1020.4
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 3061.2
$ws.Range("L113").Value = 2850
$ws.Range("M113").Value = -891.1999999999998
$ws.Range("N113").Value = -7190

$ws.Range("H122").Value = 1949.3914
$ws.Range("I122").Value = 1263.4
$ws.Range("J122").Value = 3235.625
$ws.Range("K122").Value = 3790.2
$ws.Range("L122").Value = 9706.875
$ws.Range("M122").Value = -1340.2
$ws.Range("N122").Value = -14606.875

$ws.Range("H126").Value = 3299.3333
$ws.Range("I126").Value = 1450
$ws.Range("J126").Value = 6998
$ws.Range("K126").Value = 4350
$ws.Range("L126").Value = 20994
$ws.Range("M126").Value = -1880
$ws.Range("N126").Value = -25934

$ws.Range("H136").Value = 2701.05
$ws.Range("I136").Value = 2701.05
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8103.150000000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5553.150000000001
$ws.Range("N136").ClearContents()
